$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"
$ws.Activate()
$ws.Range("I9").Select()
